# Updated capital structure database
# Applies numeric/text updates to rows 2-6 and removes row 7 (Chuou International Group)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 holds the numeric-looking text "4" (was "5"); format cell as Text first so
# Excel stores it as a string rather than re-typing it as a number, then restore
# the default style so only the underlying value (not the formatting) changes.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("B2").Style = "Normal"

# --- Row 2: refreshed metrics ---

# --- Row 2 ---
$ws.Range("D2").Value = -0.04330000000000001
$ws.Range("E2").Value = -0.03150000000000001
$ws.Range("G2").Value = 0.08759147912729771
$ws.Range("H2").Value = 0.08749355780793677
$ws.Range("I2").Value = 0.05997315660726827
$ws.Range("J2").Value = 0.04916237863417025
$ws.Range("K2").Value = 81.31
$ws.Range("L2").Value = 0.06984195155471569
$ws.Range("M2").Value = 20.62796
$ws.Range("N2").Value = 0.02063208641728346
$ws.Range("O2").Value = 0.2536952404378305
$ws.Range("P2").Value = 19.66096
$ws.Range("Q2").Value = 0.01966489297859572
$ws.Range("R2").Value = 0.2418024843192719
$ws.Range("S2").Value = 0.9670000000000005
$ws.Range("T2").Value = 0.04687812076424428
$ws.Range("U2").Value = 392.3
$ws.Range("V2").Value = 0.3923784756951391
$ws.Range("W2").Value = 0.1426665670079976
$ws.Range("X2").Value = 0.04592557063232532
$ws.Range("Y2").Value = 0.09674099637567227
$ws.Range("Z2").Value = 3.974514826613329
$ws.Range("AA2").Value = 0.2183702323550416
$ws.Range("AB2").Value = 0.04491047473445834
$ws.Range("AC2").Value = 0.1733695529513009
$ws.Range("AD2").Value = 156.7
$ws.Range("AE2").Value = 0.02125538909135998
$ws.Range("AF2").Value = 156.7212553890914
$ws.Range("AG2").Value = -235.5787446109086
$ws.Range("AH2").Value = 0.1355109166034006
$ws.Range("AI2").Value = 0.2104993702391853
$ws.Range("AJ2").Value = -0.308259869703529
$ws.Range("AK2").Value = -0.6688373884496829
$ws.Range("AL2").Value = 1.464
$ws.Range("AM2").Value = 1.36
$ws.Range("AN2").Value = 1.394686484802635
$ws.Range("AO2").Value = 47.68442622950819
$ws.Range("AP2").Value = -2.096735744834753
$ws.Range("AQ2").Value = 51.33088235294117

# --- Row 3 ---
$ws.Range("D3").Value = -0.163
$ws.Range("E3").Value = -0.07980000000000001
$ws.Range("G3").Value = 0.1702992242334688
$ws.Range("H3").Value = 0.1702992242334688
$ws.Range("I3").Value = 0.07794606575544885
$ws.Range("J3").Value = 0.07704162828360292
$ws.Range("K3").Value = 21.4
$ws.Range("L3").Value = 0.07905430365718508
$ws.Range("M3").Value = 14.5352
$ws.Range("N3").Value = 0.03635617808904452
$ws.Range("O3").Value = 0.6792149532710281
$ws.Range("P3").Value = 13.5682
$ws.Range("Q3").Value = 0.03393746873436718
$ws.Range("R3").Value = 0.6340280373831776
$ws.Range("S3").Value = 0.9670000000000005
$ws.Range("T3").Value = 0.06652815234740496
$ws.Range("U3").Value = 64.2
$ws.Range("V3").Value = 0.1605802901450726
$ws.Range("W3").Value = 0.18337617823479
$ws.Range("X3").Value = 0.05378572899931522
$ws.Range("Y3").Value = 0.1295904492354748
$ws.Range("Z3").Value = 3.682993197278912
$ws.Range("AA3").Value = 0.2837437928758001
$ws.Range("AB3").Value = 0.04432096532828639
$ws.Range("AC3").Value = 0.2394228275475138
$ws.Range("AD3").Value = 141.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 141.8
$ws.Range("AG3").Value = 77.60000000000001
$ws.Range("AH3").Value = 0.2618168389955687
$ws.Range("AI3").Value = 0.6181342632955537
$ws.Range("AJ3").Value = 0.1625471302890658
$ws.Range("AK3").Value = 0.4697336561743342
$ws.Range("AL3").Value = 1.35
$ws.Range("AM3").Value = 1.284
$ws.Range("AN3").Value = 3.244851258581236
$ws.Range("AO3").Value = 15.62962962962963
$ws.Range("AP3").Value = 1.775743707093822
$ws.Range("AQ3").Value = 16.43302180685358

# --- Row 4 ---
$ws.Range("G4").Value = 0.15435
$ws.Range("H4").Value = 0.1515
$ws.Range("I4").Value = 0.09275
$ws.Range("J4").Value = 0.06092885638297872
$ws.Range("K4").Value = 2.47
$ws.Range("L4").Value = 0.06175000000000001
$ws.Range("M4").Value = 0.97356
$ws.Range("N4").Value = 0.01187268292682927
$ws.Range("O4").Value = 0.3941538461538461
$ws.Range("P4").Value = 0.97356
$ws.Range("Q4").Value = 0.01187268292682927
$ws.Range("R4").Value = 0.3941538461538461
$ws.Range("U4").Value = 17.1
$ws.Range("V4").Value = 0.2085365853658537
$ws.Range("W4").Value = 0.08233333333333334
$ws.Range("X4").Value = 0.04500067940374067
$ws.Range("Y4").Value = 0.03733265392959267
$ws.Range("Z4").Value = 4.026170105686965
$ws.Range("AA4").Value = 0.2453099401428434
$ws.Range("AB4").Value = 0.04500067940374067
$ws.Range("AC4").Value = 0.2003092607391027
$ws.Range("AG4").Value = -17.1
$ws.Range("AJ4").Value = -0.263482280431433
$ws.Range("AK4").Value = -1.132450331125828
$ws.Range("AP4").Value = -3.143382352941177

# --- Row 5 ---
$ws.Range("B5").Value = "SBI Insurance Group Co., Ltd. (TSE:7326)"
$ws.Range("G5").Value = 0.04589468099217403
$ws.Range("H5").Value = 0.04589468099217403
$ws.Range("I5").Value = 0.04483353229871335
$ws.Range("J5").Value = 0.04270939677969192
$ws.Range("K5").Value = 50.8
$ws.Range("L5").Value = 0.06738294203475262
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 286
$ws.Range("V5").Value = 0.8988057825267128
$ws.Range("W5").Value = 0.1290978398983481
$ws.Range("X5").Value = 0.04500067940374067
$ws.Range("Y5").Value = 0.08409716049460747
$ws.Range("Z5").Value = 4.482164090368609
$ws.Range("AA5").Value = 0.1914305245672398
$ws.Range("AB5").Value = 0.04500067940374067
$ws.Range("AC5").Value = 0.1464298451634992
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -286
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = -8.881987577639755
$ws.Range("AK5").Value = -2.220496894409937
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = 0
$ws.Range("AP5").Value = -5.945945945945946

# --- Row 6 ---
$ws.Range("B6").Value = "Advance Create Co., Ltd. (TSE:8798)"
$ws.Range("G6").Value = 0.1516064257028112
$ws.Range("H6").Value = 0.1516064257028112
$ws.Range("I6").Value = 0.1125577201022262
$ws.Range("J6").Value = 0.07665469348500327
$ws.Range("K6").Value = 6.64
$ws.Range("L6").Value = 0.06666666666666667
$ws.Range("M6").Value = 5.1192
$ws.Range("N6").Value = 0.02562162162162162
$ws.Range("O6").Value = 0.7709638554216868
$ws.Range("P6").Value = 5.1192
$ws.Range("Q6").Value = 0.02562162162162162
$ws.Range("R6").Value = 0.7709638554216868
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 25
$ws.Range("V6").Value = 0.1251251251251251
$ws.Range("W6").Value = 0.1562352941176471
$ws.Range("X6").Value = 0.04685046186090998
$ws.Range("Y6").Value = 0.1093848322567371
$ws.Range("Z6").Value = 2.412717323182944
$ws.Range("AA6").Value = 0.1849461068745461
$ws.Range("AB6").Value = 0.04482027006517601
$ws.Range("AC6").Value = 0.1401258368093701
$ws.Range("AD6").Value = 14.9
$ws.Range("AE6").Value = 0.02125538909135998
$ws.Range("AF6").Value = 14.92125538909136
$ws.Range("AG6").Value = -10.07874461090864
$ws.Range("AH6").Value = 0.0694912823700332
$ws.Range("AI6").Value = 0.2190396419423707
$ws.Range("AJ6").Value = -0.05312396120423389
$ws.Range("AK6").Value = -0.2337303151303503
$ws.Range("AL6").Value = 0.114
$ws.Range("AM6").Value = 0.07600000000000001
$ws.Range("AN6").Value = 0.9857757194839564
$ws.Range("AP6").Value = -0.6668041423029203

# --- Row 6 gained historical growth + interest-coverage figures ---
$ws.Range("D6").Value = 0.0764
$ws.Range("E6").Value = 0.0168
$ws.Range("AO6").Value = 98.24561403508771
$ws.Range("AQ6").Value = 147.3684210526315

# --- Row 5 lost historical growth + buyback/coverage figures (now blank) ---
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AQ5").ClearContents()

# --- Remove the old row 7 (Chuou International Group Co., Ltd.) entirely ---
$ws.Rows(7).EntireRow.Delete()
